$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H106").Value = 3554.4614
$ws.Range("I106").Value = 2912.111
$ws.Range("K106").Value = 2912.111
$ws.Range("M106").Value = -2281.111

$ws.Range("H129").Value = 2549.7585
$ws.Range("I129").Value = 714.9
$ws.Range("J129").Value = 3515.4736
$ws.Range("K129").Value = 2144.7
$ws.Range("L129").Value = 10546.4208
$ws.Range("M129").Value = 2855.3
$ws.Range("N129").Value = -20546.4208

$ws.Range("H138").Value = 4886.2593
$ws.Range("I138").Value = 3445.5334
$ws.Range("J138").Value = 5440.385
$ws.Range("K138").Value = 10336.6002
$ws.Range("L138").Value = 16321.155
$ws.Range("M138").Value = -5196.600199999999
$ws.Range("N138").Value = -26601.155

$ws.Range("H141").Value = 1540
$ws.Range("I141").Value = 1540
$ws.Range("K141").Value = 4620
$ws.Range("M141").Value = 560

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 6098573.5
$ws.Range("I32").Value = 6938032
$ws.Range("K32").Value = 6938032
$ws.Range("M32").Value = -6937745

$ws.Range("H61").Value = 2758.724
$ws.Range("I61").Value = 1900.125
$ws.Range("K61").Value = 1900.125
$ws.Range("M61").Value = -1688.125

$ws.Range("H74").Value = 4872.026
$ws.Range("I74").Value = 3930.3
$ws.Range("K74").Value = 3930.3
$ws.Range("M74").Value = -3056.3

$ws.Range("H77").Value = 4872.026
$ws.Range("I77").Value = 3930.3
$ws.Range("K77").Value = 19651.5
$ws.Range("M77").Value = -15283.5

$ws.Range("H95").Value = 45000
$ws.Range("J95").Value = 45000
$ws.Range("L95").Value = 45000
$ws.Range("N95").Value = -50492

$ws.Range("H132").Value = 4010.9844
$ws.Range("I132").Value = 2132.7144
$ws.Range("J132").Value = 10146.667
$ws.Range("K132").Value = 6398.1432
$ws.Range("L132").Value = 30440.001
$ws.Range("M132").Value = -3868.1432
$ws.Range("N132").Value = -35500.001

$ws.Range("H136").Value = 2758.724
$ws.Range("I136").Value = 1900.125
$ws.Range("K136").Value = 5700.375
$ws.Range("M136").Value = -3150.375

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H95").Value = 14668
$ws.Range("J95").Value = 14668
$ws.Range("L95").Value = 14668
$ws.Range("N95").Value = -20160

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H68").Value = 2969.75
$ws.Range("J68").Value = 3183.3928
$ws.Range("L68").Value = 9550.178400000001
$ws.Range("N68").Value = -11172.1784

$ws.Range("H71").Value = 2969.75
$ws.Range("J71").Value = 3183.3928
$ws.Range("L71").Value = 28650.5352
$ws.Range("N71").Value = -36762.5352

$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0

$ws.Range("H107").Value = 7541.091
$ws.Range("J107").Value = 12498
$ws.Range("L107").Value = 37494
$ws.Range("N107").Value = -41334

$ws.Range("H136").Value = 4012
$ws.Range("I136").Value = 4012
$ws.Range("K136").Value = 12036
$ws.Range("M136").Value = -6936

$ws.Range("H137").Value = 1599.6666
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0

$ws.Range("H140").Value = 50000890
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0

$ws.Range("N104").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("N140").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 9743.467000000001
$ws.Range("I70").Value = 9085.299999999999
$ws.Range("J70").Value = 11059.8
$ws.Range("K70").Value = 9085.299999999999
$ws.Range("L70").Value = 11059.8
$ws.Range("M70").Value = -8815.299999999999
$ws.Range("N70").Value = -11599.8

$ws.Range("H73").Value = 9743.467000000001
$ws.Range("I73").Value = 9085.299999999999
$ws.Range("J73").Value = 11059.8
$ws.Range("K73").Value = 9085.299999999999
$ws.Range("L73").Value = 11059.8
$ws.Range("M73").Value = -8149.299999999999
$ws.Range("N73").Value = -12931.8

$ws.Range("H122").Value = 7191.7856
$ws.Range("I122").Value = 8918.6
$ws.Range("J122").Value = 2874.75
$ws.Range("K122").Value = 26755.8
$ws.Range("L122").Value = 8624.25
$ws.Range("M122").Value = -24305.8
$ws.Range("N122").Value = -13524.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H40").Value = 9334.666999999999
$ws.Range("I40").Value = 6502
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 6502
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = -6366
$ws.Range("N40").Value = -15272

$ws.Range("H61").Value = 5157.5
$ws.Range("I61").Value = 5071.364
$ws.Range("K61").Value = 5071.364
$ws.Range("M61").Value = -4869.364

$ws.Range("H113").Value = 5157.5
$ws.Range("I113").Value = 5071.364
$ws.Range("K113").Value = 5071.364
$ws.Range("M113").Value = -2901.364

$ws.Range("H122").Value = 250002030
$ws.Range("I122").Value = 500000000
$ws.Range("J122").Value = 4050
$ws.Range("K122").Value = 1500000000
$ws.Range("L122").Value = 12150
$ws.Range("M122").Value = -1499997550
$ws.Range("N122").Value = -17050

$ws.Range("H132").Value = 4984.5
$ws.Range("I132").Value = 3181.5625
$ws.Range("K132").Value = 9544.6875
$ws.Range("M132").Value = -7014.6875

$ws.Range("H136").Value = 5264.921
$ws.Range("I136").Value = 2609.158
$ws.Range("J136").Value = 7920.684
$ws.Range("K136").Value = 7827.474
$ws.Range("L136").Value = 23762.052
$ws.Range("M136").Value = -5277.474
$ws.Range("N136").Value = -28862.052

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H2").Value = 3249999.8
$ws.Range("I2").Value = 5500000
$ws.Range("K2").Value = 5500000
$ws.Range("M2").Value = -5499888

$ws.Range("H95").Value = 42053.832
$ws.Range("J95").Value = 42053.832
$ws.Range("L95").Value = 42053.832
$ws.Range("N95").Value = -47545.832

$ws.Range("H103").Value = 56051.855
$ws.Range("J103").Value = 56051.855
$ws.Range("L103").Value = 56051.855
$ws.Range("N103").Value = -58395.855

$ws.Range("H122").Value = 3724.6428
$ws.Range("I122").Value = 3467.75
$ws.Range("K122").Value = 10403.25
$ws.Range("M122").Value = -7953.25

$ws.Range("H126").Value = 3742913.8
$ws.Range("I126").Value = 4592778.5
$ws.Range("K126").Value = 13778335.5
$ws.Range("M126").Value = -13775865.5
